# Generate Report for Handoff
# Updates the localization-status workbook so that "b.md" is reported as
# "Ready for handoff" (instead of "Handed back: in sync with en-US") on
# the Overview sheet as well as the zh-cn / de-de detail sheets, and adds
# the corresponding handoff-file / datetime / error-detail information.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (row 3 = b.md) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-26 00:37:12"

# ---- zh-cn sheet (row 3 = b.md) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-26 00:37:07"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc90e69b0c4083b0693ebb499399e1660847f6c0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07e27dc9a358853099b9fdf1ee9433730e8be43b/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---- de-de sheet (row 3 = b.md) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-26 00:37:12"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc90e69b0c4083b0693ebb499399e1660847f6c0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07e27dc9a358853099b9fdf1ee9433730e8be43b/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
